$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 195, shifting existing rows 195:258 down to 196:259.
$ws.Rows("195:195").Insert()

# Populate the newly inserted row with the new data record.
$ws.Range("A195").Value = 10
$ws.Range("B195").Value = "Vega Modelo de Temuco"
$ws.Range("C195").Value = "La Araucanía"
$ws.Range("D195").Value = 44588
$ws.Range("E195").Value = 9
$ws.Range("F195").Value = 100114013
$ws.Range("G195").Value = "Zanahoria"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 125
$ws.Range("K195").Value = 9000
$ws.Range("L195").Value = 9000
$ws.Range("M195").Value = 9000
$ws.Range("N195").Value = "$/saco 20 kilos"
$ws.Range("O195").Value = "Región del Maule"
$ws.Range("P195").Value = 450
$ws.Range("Q195").Value = 20
$ws.Range("R195").Value = "Hortaliza"
